# Re-sync DATOS_FINANCIEROS against the refreshed Supabase export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rows that were re-ordered (E/F/G swapped between the two sibling rows) ---
# Rows 24-25
    $ws.Cells.Item(24, 5).Value = 101
    $ws.Cells.Item(24, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(24, 7).Value = 12234056
    $ws.Cells.Item(25, 5).Value = 104
    $ws.Cells.Item(25, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(25, 7).Value = 0
# Rows 34-35
    $ws.Cells.Item(34, 5).Value = 101
    $ws.Cells.Item(34, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(35, 5).Value = 104
    $ws.Cells.Item(35, 6).Value = 'NUEVOS PROYECTOS'
# Rows 71-72
    $ws.Cells.Item(71, 5).Value = 10203
    $ws.Cells.Item(71, 6).Value = 'SISTEMAS'
    $ws.Cells.Item(71, 7).Value = 321116
    $ws.Cells.Item(72, 5).Value = 0
    $ws.Cells.Item(72, 6).Value = 'VARIOS'
    $ws.Cells.Item(72, 7).Value = 117563
# Rows 148-149
    $ws.Cells.Item(148, 5).Value = 104
    $ws.Cells.Item(148, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(149, 5).Value = 101
    $ws.Cells.Item(149, 6).Value = 'V.P. PROYECTOS'
# Rows 195-196
    $ws.Cells.Item(195, 5).Value = 101
    $ws.Cells.Item(195, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(195, 7).Value = 11552916
    $ws.Cells.Item(196, 5).Value = 104
    $ws.Cells.Item(196, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(196, 7).Value = 0
# Rows 250-251
    $ws.Cells.Item(250, 5).Value = 104
    $ws.Cells.Item(250, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(250, 7).Value = 0
    $ws.Cells.Item(251, 5).Value = 101
    $ws.Cells.Item(251, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(251, 7).Value = 12548548
# Rows 257-258
    $ws.Cells.Item(257, 5).Value = 101
    $ws.Cells.Item(257, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(258, 5).Value = 104
    $ws.Cells.Item(258, 6).Value = 'NUEVOS PROYECTOS'
# Rows 311-312
    $ws.Cells.Item(311, 5).Value = 101
    $ws.Cells.Item(311, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(312, 5).Value = 104
    $ws.Cells.Item(312, 6).Value = 'NUEVOS PROYECTOS'
# Rows 358-359
    $ws.Cells.Item(358, 5).Value = 101
    $ws.Cells.Item(358, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(358, 7).Value = 10916278
    $ws.Cells.Item(359, 5).Value = 104
    $ws.Cells.Item(359, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(359, 7).Value = 0
# Rows 365-366
    $ws.Cells.Item(365, 5).Value = 101
    $ws.Cells.Item(365, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(366, 5).Value = 104
    $ws.Cells.Item(366, 6).Value = 'NUEVOS PROYECTOS'
# Rows 408-409
    $ws.Cells.Item(408, 5).Value = 101
    $ws.Cells.Item(408, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(408, 7).Value = 12548281
    $ws.Cells.Item(409, 5).Value = 104
    $ws.Cells.Item(409, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(409, 7).Value = 0
# Rows 462-463
    $ws.Cells.Item(462, 5).Value = 101
    $ws.Cells.Item(462, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(462, 7).Value = 12548281
    $ws.Cells.Item(463, 5).Value = 104
    $ws.Cells.Item(463, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(463, 7).Value = 0
# Rows 468-469
    $ws.Cells.Item(468, 5).Value = 101
    $ws.Cells.Item(468, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(469, 5).Value = 104
    $ws.Cells.Item(469, 6).Value = 'NUEVOS PROYECTOS'
# Rows 516-517
    $ws.Cells.Item(516, 5).Value = 101
    $ws.Cells.Item(516, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(516, 7).Value = 10238677.4
    $ws.Cells.Item(517, 5).Value = 104
    $ws.Cells.Item(517, 6).Value = 'NUEVOS PROYECTOS'
    $ws.Cells.Item(517, 7).Value = 0
# Rows 522-523
    $ws.Cells.Item(522, 5).Value = 101
    $ws.Cells.Item(522, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(523, 5).Value = 104
    $ws.Cells.Item(523, 6).Value = 'NUEVOS PROYECTOS'

# --- May-2025 (and new Jun-2025) account rows re-pulled from source: shift in
#     new accounts (4175, 4250, 4295, 5150, 5215, 5405...) and append trailing rows ---
# Row 876
    $ws.Cells.Item(876, 3).Value = 4155
    $ws.Cells.Item(876, 4).Value = 'ACTIVIDADES INMOBILIARIAS; EMPRESARIALES Y DE ALQUILER'
    $ws.Cells.Item(876, 7).Value = 0
    $ws.Cells.Item(876, 8).Value = 1659058403
# Row 877
    $ws.Cells.Item(877, 3).Value = 4175
    $ws.Cells.Item(877, 4).Value = 'DEVOLUCIONES; REBAJAS Y DESCUENTOS EN VENTAS (DB)'
    $ws.Cells.Item(877, 7).Value = 16577200
    $ws.Cells.Item(877, 8).Value = 0
# Row 878
    $ws.Cells.Item(878, 3).Value = 4210
    $ws.Cells.Item(878, 4).Value = 'FINANCIEROS'
    $ws.Cells.Item(878, 5).Value = 0
    $ws.Cells.Item(878, 6).Value = 'VARIOS'
    $ws.Cells.Item(878, 7).Value = 1280429.26
    $ws.Cells.Item(878, 8).Value = 18490297.350000001
# Row 879
    $ws.Cells.Item(879, 3).Value = 4250
    $ws.Cells.Item(879, 4).Value = 'RECUPERACIONES'
    $ws.Cells.Item(879, 5).Value = 0
    $ws.Cells.Item(879, 6).Value = 'VARIOS'
    $ws.Cells.Item(879, 7).Value = 0
    $ws.Cells.Item(879, 8).Value = 2072017
# Row 880
    $ws.Cells.Item(880, 3).Value = 4295
    $ws.Cells.Item(880, 4).Value = 'DIVERSOS'
    $ws.Cells.Item(880, 5).Value = 0
    $ws.Cells.Item(880, 6).Value = 'VARIOS'
    $ws.Cells.Item(880, 7).Value = 0
    $ws.Cells.Item(880, 8).Value = 2230.1999999999998
# Row 881
    $ws.Cells.Item(881, 3).Value = 5105
    $ws.Cells.Item(881, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(881, 7).Value = 8854827
    $ws.Cells.Item(881, 8).Value = 4075179
# Row 882
    $ws.Cells.Item(882, 3).Value = 5105
    $ws.Cells.Item(882, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(882, 5).Value = 100
    $ws.Cells.Item(882, 6).Value = 'V.P. ESTRATEGICA'
    $ws.Cells.Item(882, 7).Value = 43902068
# Row 883
    $ws.Cells.Item(883, 3).Value = 5105
    $ws.Cells.Item(883, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(883, 5).Value = 101
    $ws.Cells.Item(883, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(883, 7).Value = 18186097
# Row 884
    $ws.Cells.Item(884, 3).Value = 5105
    $ws.Cells.Item(884, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(884, 5).Value = 102
    $ws.Cells.Item(884, 6).Value = 'ADMINISTRACION'
    $ws.Cells.Item(884, 7).Value = 49120753.039999999
    $ws.Cells.Item(884, 8).Value = 0
# Row 885
    $ws.Cells.Item(885, 3).Value = 5110
    $ws.Cells.Item(885, 4).Value = 'HONORARIOS'
    $ws.Cells.Item(885, 7).Value = 8661200
# Row 886
    $ws.Cells.Item(886, 3).Value = 5130
    $ws.Cells.Item(886, 4).Value = 'SEGURO'
    $ws.Cells.Item(886, 7).Value = 45460
# Row 887
    $ws.Cells.Item(887, 3).Value = 5130
    $ws.Cells.Item(887, 4).Value = 'SEGURO'
    $ws.Cells.Item(887, 5).Value = 102
    $ws.Cells.Item(887, 6).Value = 'ADMINISTRACION'
    $ws.Cells.Item(887, 7).Value = 31540
# Row 888
    $ws.Cells.Item(888, 3).Value = 5135
    $ws.Cells.Item(888, 4).Value = 'SERVICIOS'
    $ws.Cells.Item(888, 7).Value = 12031128
    $ws.Cells.Item(888, 8).Value = 1671135.7
# Row 889
    $ws.Cells.Item(889, 3).Value = 5140
    $ws.Cells.Item(889, 4).Value = 'GASTOS LEGALES'
    $ws.Cells.Item(889, 7).Value = 508887
# Row 890
    $ws.Cells.Item(890, 3).Value = 5145
    $ws.Cells.Item(890, 4).Value = 'MANTENIMIENTO Y REPARACIONES'
    $ws.Cells.Item(890, 7).Value = 3476065
    $ws.Cells.Item(890, 8).Value = 0
# Row 891
    $ws.Cells.Item(891, 3).Value = 5150
    $ws.Cells.Item(891, 4).Value = 'ADECUACION E INSTALACION'
    $ws.Cells.Item(891, 5).Value = 10507008
    $ws.Cells.Item(891, 6).Value = 'MANTENIMENTOS PLANTA'
    $ws.Cells.Item(891, 7).Value = 38915.769999999997
# Row 892
    $ws.Cells.Item(892, 3).Value = 5155
    $ws.Cells.Item(892, 4).Value = 'GASTOS DE VIAJE'
    $ws.Cells.Item(892, 5).Value = 0
    $ws.Cells.Item(892, 6).Value = 'VARIOS'
    $ws.Cells.Item(892, 7).Value = 7229826
# Row 893
    $ws.Cells.Item(893, 3).Value = 5160
    $ws.Cells.Item(893, 4).Value = 'DEPRECIACIONES'
    $ws.Cells.Item(893, 5).Value = 0
    $ws.Cells.Item(893, 6).Value = 'VARIOS'
    $ws.Cells.Item(893, 7).Value = 7259515.8700000001
# Row 894
    $ws.Cells.Item(894, 3).Value = 5165
    $ws.Cells.Item(894, 4).Value = 'AMORTIZACIONES'
    $ws.Cells.Item(894, 5).Value = 0
    $ws.Cells.Item(894, 6).Value = 'VARIOS'
    $ws.Cells.Item(894, 7).Value = 3423821.48
# Row 895
    $ws.Cells.Item(895, 3).Value = 5195
    $ws.Cells.Item(895, 4).Value = 'DIVERSOS'
    $ws.Cells.Item(895, 5).Value = 0
    $ws.Cells.Item(895, 6).Value = 'VARIOS'
    $ws.Cells.Item(895, 7).Value = 14304598
    $ws.Cells.Item(895, 8).Value = 242547
# Row 896
    $ws.Cells.Item(896, 5).Value = 0
    $ws.Cells.Item(896, 6).Value = 'VARIOS'
    $ws.Cells.Item(896, 7).Value = 10297258
# Row 897
    $ws.Cells.Item(897, 5).Value = 105
    $ws.Cells.Item(897, 6).Value = 'DISEÑO'
    $ws.Cells.Item(897, 7).Value = 31790793
# Row 898
    $ws.Cells.Item(898, 3).Value = 5205
    $ws.Cells.Item(898, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(898, 5).Value = 100
    $ws.Cells.Item(898, 6).Value = 'V.P. ESTRATEGICA'
    $ws.Cells.Item(898, 7).Value = 0
    $ws.Cells.Item(898, 8).Value = 0
# Row 899
    $ws.Cells.Item(899, 3).Value = 5205
    $ws.Cells.Item(899, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(899, 5).Value = 102
    $ws.Cells.Item(899, 6).Value = 'ADMINISTRACION'
    $ws.Cells.Item(899, 7).Value = 8249385.4500000002
# Row 900
    $ws.Cells.Item(900, 3).Value = 5205
    $ws.Cells.Item(900, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(900, 5).Value = 101
    $ws.Cells.Item(900, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(900, 7).Value = 0
# Row 901
    $ws.Cells.Item(901, 3).Value = 5205
    $ws.Cells.Item(901, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(901, 7).Value = 2264161.2000000002
# Row 902
    $ws.Cells.Item(902, 3).Value = 5205
    $ws.Cells.Item(902, 4).Value = 'GASTOS DE PERSONAL'
    $ws.Cells.Item(902, 5).Value = 106
    $ws.Cells.Item(902, 6).Value = 'COMERCIAL'
    $ws.Cells.Item(902, 7).Value = 41602585.990000002
# Row 903
    $ws.Cells.Item(903, 3).Value = 5210
    $ws.Cells.Item(903, 4).Value = 'HONORARIOS'
    $ws.Cells.Item(903, 5).Value = 0
    $ws.Cells.Item(903, 6).Value = 'VARIOS'
    $ws.Cells.Item(903, 7).Value = 17083600
    $ws.Cells.Item(903, 8).Value = 749550
# Row 904
    $ws.Cells.Item(904, 3).Value = 5215
    $ws.Cells.Item(904, 4).Value = 'IMPUESTOS'
    $ws.Cells.Item(904, 5).Value = 0
    $ws.Cells.Item(904, 6).Value = 'VARIOS'
    $ws.Cells.Item(904, 7).Value = 10066126
# Row 905
    $ws.Cells.Item(905, 3).Value = 5230
    $ws.Cells.Item(905, 4).Value = 'SEGURO'
    $ws.Cells.Item(905, 5).Value = 105
    $ws.Cells.Item(905, 6).Value = 'DISEÑO'
    $ws.Cells.Item(905, 7).Value = 10880
# Row 906
    $ws.Cells.Item(906, 3).Value = 5230
    $ws.Cells.Item(906, 4).Value = 'SEGURO'
    $ws.Cells.Item(906, 5).Value = 100
    $ws.Cells.Item(906, 6).Value = 'V.P. ESTRATEGICA'
    $ws.Cells.Item(906, 7).Value = 5440
# Row 907
    $ws.Cells.Item(907, 3).Value = 5230
    $ws.Cells.Item(907, 4).Value = 'SEGURO'
    $ws.Cells.Item(907, 5).Value = 107
    $ws.Cells.Item(907, 6).Value = 'PLANTA PRODUCCION'
    $ws.Cells.Item(907, 7).Value = 5440
# Row 908
    $ws.Cells.Item(908, 3).Value = 5230
    $ws.Cells.Item(908, 4).Value = 'SEGURO'
    $ws.Cells.Item(908, 5).Value = 101
    $ws.Cells.Item(908, 6).Value = 'V.P. PROYECTOS'
    $ws.Cells.Item(908, 7).Value = 5440
    $ws.Cells.Item(908, 8).Value = 0
# Row 909
    $ws.Cells.Item(909, 3).Value = 5230
    $ws.Cells.Item(909, 4).Value = 'SEGURO'
    $ws.Cells.Item(909, 5).Value = 102
    $ws.Cells.Item(909, 6).Value = 'ADMINISTRACION'
    $ws.Cells.Item(909, 7).Value = 10880
# Row 910
    $ws.Cells.Item(910, 3).Value = 5235
    $ws.Cells.Item(910, 4).Value = 'SERVICIOS'
    $ws.Cells.Item(910, 5).Value = 102
    $ws.Cells.Item(910, 6).Value = 'ADMINISTRACION'
    $ws.Cells.Item(910, 7).Value = 33000
# Row 911
    $ws.Cells.Item(911, 3).Value = 5235
    $ws.Cells.Item(911, 4).Value = 'SERVICIOS'
    $ws.Cells.Item(911, 5).Value = 107
    $ws.Cells.Item(911, 6).Value = 'PLANTA PRODUCCION'
    $ws.Cells.Item(911, 7).Value = 33000
# Row 912
    $ws.Cells.Item(912, 3).Value = 5235
    $ws.Cells.Item(912, 4).Value = 'SERVICIOS'
    $ws.Cells.Item(912, 5).Value = 106
    $ws.Cells.Item(912, 6).Value = 'COMERCIAL'
    $ws.Cells.Item(912, 7).Value = 30000
# Row 913
    $ws.Cells.Item(913, 3).Value = 5235
    $ws.Cells.Item(913, 4).Value = 'SERVICIOS'
    $ws.Cells.Item(913, 5).Value = 105
    $ws.Cells.Item(913, 6).Value = 'DISEÑO'
    $ws.Cells.Item(913, 7).Value = 33000
# Row 914
    $ws.Cells.Item(914, 3).Value = 5235
    $ws.Cells.Item(914, 4).Value = 'SERVICIOS'
    $ws.Cells.Item(914, 7).Value = 5659185
    $ws.Cells.Item(914, 8).Value = 1671135.7
# Row 915
    $ws.Cells.Item(915, 3).Value = 5250
    $ws.Cells.Item(915, 4).Value = 'ADECUACION E INSTALACION'
    $ws.Cells.Item(915, 5).Value = 10507008
    $ws.Cells.Item(915, 6).Value = 'MANTENIMENTOS PLANTA'
    $ws.Cells.Item(915, 7).Value = 38915.769999999997
# Row 916
    $ws.Cells.Item(916, 3).Value = 5255
    $ws.Cells.Item(916, 4).Value = 'GASTOS DE VIAJE'
    $ws.Cells.Item(916, 7).Value = 190739
    $ws.Cells.Item(916, 8).Value = 0
# Row 917
    $ws.Cells.Item(917, 3).Value = 5260
    $ws.Cells.Item(917, 4).Value = 'DEPRECIACIONES'
    $ws.Cells.Item(917, 7).Value = 2216177.29
    $ws.Cells.Item(917, 8).Value = 0
# Row 918
    $ws.Cells.Item(918, 3).Value = 5265
    $ws.Cells.Item(918, 4).Value = 'AMORTIZACIONES'
    $ws.Cells.Item(918, 7).Value = 4710220.1100000003
    $ws.Cells.Item(918, 8).Value = 0
# Row 919
    $ws.Cells.Item(919, 3).Value = 5295
    $ws.Cells.Item(919, 4).Value = 'DIVERSOS'
    $ws.Cells.Item(919, 7).Value = 7626408
# Row 920
    $ws.Cells.Item(920, 1).Value = 2025
    $ws.Cells.Item(920, 2).Value = 5
    $ws.Cells.Item(920, 3).Value = 5305
    $ws.Cells.Item(920, 4).Value = 'FINANCIEROS'
    $ws.Cells.Item(920, 5).Value = 0
    $ws.Cells.Item(920, 6).Value = 'VARIOS'
    $ws.Cells.Item(920, 7).Value = 25497977.789999999
    $ws.Cells.Item(920, 8).Value = 0
# Row 921
    $ws.Cells.Item(921, 1).Value = 2025
    $ws.Cells.Item(921, 2).Value = 5
    $ws.Cells.Item(921, 3).Value = 5315
    $ws.Cells.Item(921, 4).Value = 'GASTOS EXTRAORDINARIOS'
    $ws.Cells.Item(921, 5).Value = 0
    $ws.Cells.Item(921, 6).Value = 'VARIOS'
    $ws.Cells.Item(921, 7).Value = 88681067.030000001
    $ws.Cells.Item(921, 8).Value = 0
# Row 922
    $ws.Cells.Item(922, 1).Value = 2025
    $ws.Cells.Item(922, 2).Value = 5
    $ws.Cells.Item(922, 3).Value = 5395
    $ws.Cells.Item(922, 4).Value = 'GASTOS DIVERSOS'
    $ws.Cells.Item(922, 5).Value = 0
    $ws.Cells.Item(922, 6).Value = 'VARIOS'
    $ws.Cells.Item(922, 7).Value = 2265.38
    $ws.Cells.Item(922, 8).Value = 0
# Row 923
    $ws.Cells.Item(923, 1).Value = 2025
    $ws.Cells.Item(923, 2).Value = 5
    $ws.Cells.Item(923, 3).Value = 5405
    $ws.Cells.Item(923, 4).Value = 'IMPUESTO DE RENTA Y COMPLEMENTARIOS'
    $ws.Cells.Item(923, 5).Value = 0
    $ws.Cells.Item(923, 6).Value = 'VARIOS'
    $ws.Cells.Item(923, 7).Value = 41727074
    $ws.Cells.Item(923, 8).Value = 0
# Row 924
    $ws.Cells.Item(924, 1).Value = 2025
    $ws.Cells.Item(924, 2).Value = 5
    $ws.Cells.Item(924, 3).Value = 615555
    $ws.Cells.Item(924, 4).Value = 'COSTO MATERIALES'
    $ws.Cells.Item(924, 5).Value = 0
    $ws.Cells.Item(924, 6).Value = 'VARIOS'
    $ws.Cells.Item(924, 7).Value = 1265462273.3499999
    $ws.Cells.Item(924, 8).Value = 130728360.03
# Row 925
    $ws.Cells.Item(925, 1).Value = 2025
    $ws.Cells.Item(925, 2).Value = 5
    $ws.Cells.Item(925, 3).Value = 615555
    $ws.Cells.Item(925, 4).Value = 'COSTO MATERIALES'
    $ws.Cells.Item(925, 5).Value = 10507010
    $ws.Cells.Item(925, 6).Value = 'PRODUCCION'
    $ws.Cells.Item(925, 7).Value = 15250180.060000001
    $ws.Cells.Item(925, 8).Value = 0
# Row 926
    $ws.Cells.Item(926, 1).Value = 2025
    $ws.Cells.Item(926, 2).Value = 5
    $ws.Cells.Item(926, 3).Value = 615556
    $ws.Cells.Item(926, 4).Value = 'PUBLICIDAD - SERVICIOS'
    $ws.Cells.Item(926, 5).Value = 0
    $ws.Cells.Item(926, 6).Value = 'VARIOS'
    $ws.Cells.Item(926, 7).Value = 12227220
    $ws.Cells.Item(926, 8).Value = 40410477.399999999
# Row 927
    $ws.Cells.Item(927, 1).Value = 2025
    $ws.Cells.Item(927, 2).Value = 5
    $ws.Cells.Item(927, 3).Value = 615558
    $ws.Cells.Item(927, 4).Value = 'COSTO MANO DE OBRA DIRECTA'
    $ws.Cells.Item(927, 5).Value = 0
    $ws.Cells.Item(927, 6).Value = 'VARIOS'
    $ws.Cells.Item(927, 7).Value = 135224635.25
    $ws.Cells.Item(927, 8).Value = 580905.06999999995
# Row 928
    $ws.Cells.Item(928, 1).Value = 2025
    $ws.Cells.Item(928, 2).Value = 5
    $ws.Cells.Item(928, 3).Value = 615559
    $ws.Cells.Item(928, 4).Value = 'COSTO MANO DE OBRA INDIRECTA'
    $ws.Cells.Item(928, 5).Value = 0
    $ws.Cells.Item(928, 6).Value = 'VARIOS'
    $ws.Cells.Item(928, 7).Value = 43615469.460000001
    $ws.Cells.Item(928, 8).Value = 0
# Row 929
    $ws.Cells.Item(929, 1).Value = 2025
    $ws.Cells.Item(929, 2).Value = 5
    $ws.Cells.Item(929, 3).Value = 615560
    $ws.Cells.Item(929, 4).Value = 'COSTO CIF'
    $ws.Cells.Item(929, 5).Value = 0
    $ws.Cells.Item(929, 6).Value = 'VARIOS'
    $ws.Cells.Item(929, 7).Value = 190430566.24000001
    $ws.Cells.Item(929, 8).Value = 372454809.5
# Row 930
    $ws.Cells.Item(930, 1).Value = 2025
    $ws.Cells.Item(930, 2).Value = 5
    $ws.Cells.Item(930, 3).Value = 615561
    $ws.Cells.Item(930, 4).Value = 'MANO DE OBRA DIRECTA - TIEMPOS MUERTOS'
    $ws.Cells.Item(930, 5).Value = 0
    $ws.Cells.Item(930, 6).Value = 'VARIOS'
    $ws.Cells.Item(930, 7).Value = 0
    $ws.Cells.Item(930, 8).Value = 0
# Row 931
    $ws.Cells.Item(931, 1).Value = 2025
    $ws.Cells.Item(931, 2).Value = 6
    $ws.Cells.Item(931, 3).Value = 4210
    $ws.Cells.Item(931, 4).Value = 'FINANCIEROS'
    $ws.Cells.Item(931, 5).Value = 0
    $ws.Cells.Item(931, 6).Value = 'VARIOS'
    $ws.Cells.Item(931, 7).Value = 0
    $ws.Cells.Item(931, 8).Value = 0

# --- Restore the author's last on-screen selection (whole second row) ---
$ws.Rows("2").Select()
